# Live-coding update: add a new "Heading 1" section right after the
# date line, and anchor an R-Markdown-style bookmark ("heading-1") to
# it that spans from the new heading through to the end of the
# document (mirroring the existing "r-markdown" bookmark already
# present in the file, which spans its whole section the same way).

$d = $word.ActiveDocument

# Locate the "Date" paragraph (holds "10/12/2021") that the new
# heading paragraph must follow.
$searchRange = $d.Content
$searchRange.Find.Execute("10/12/2021", $true, $false, $false, $false,
                           $false, $true, 1, $false, "", 0)
$datePara = $searchRange.Paragraphs(1)

# Insert a new empty paragraph right after it, then give it the
# Heading1 style and the "Heading 1" text.
$datePara.Range.InsertParagraphAfter()
$headingPara = $datePara.Next()
$headingPara.Range.Text = "Heading 1"
$headingPara.Style = "Heading1"

# Bookmark the new heading paragraph through the end of the document,
# the same way the pre-existing "r-markdown" bookmark wraps its
# section.
$bmRange = $d.Range($headingPara.Range.Start, $d.Content.End)
$d.Bookmarks.Add("heading-1", $bmRange)
